$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (2-9) with new forecast-error values ---
$ws.Range("B2").Value = -0.02074140170511275
$ws.Range("C2").Value = 1.447200151776814
$ws.Range("D2").Value = 4.044434787740455
$ws.Range("E2").Value = 2.011078016323697
$ws.Range("F2").Value = 2.030981209572279
$ws.Range("G2").Value = 51

$ws.Range("B3").Value = 0.3375535848015695
$ws.Range("C3").Value = 1.320417681676975
$ws.Range("D3").Value = 3.736644118779018
$ws.Range("E3").Value = 1.93304012342709
$ws.Range("F3").Value = 1.922663356868149
$ws.Range("G3").Value = 50

$ws.Range("B4").Value = 0.09841972257903492
$ws.Range("C4").Value = 1.321951009708453
$ws.Range("D4").Value = 3.706477144602156
$ws.Range("E4").Value = 1.925221323537155
$ws.Range("F4").Value = 1.942628934146926
$ws.Range("G4").Value = 49

$ws.Range("B5").Value = 0.3286888119229458
$ws.Range("C5").Value = 1.355285765321783
$ws.Range("D5").Value = 3.906105706609585
$ws.Range("E5").Value = 1.976387033606926
$ws.Range("F5").Value = 1.969487079966377
$ws.Range("G5").Value = 48

$ws.Range("B6").Value = 0.1758242421875162
$ws.Range("C6").Value = 1.38590038373524
$ws.Range("D6").Value = 4.013282287792838
$ws.Range("E6").Value = 2.003317819965878
$ws.Range("F6").Value = 2.01716170490405
$ws.Range("G6").Value = 47

$ws.Range("B7").Value = 0.3886596252776647
$ws.Range("C7").Value = 1.341425932669102
$ws.Range("D7").Value = 3.575002999449235
$ws.Range("E7").Value = 1.890767833301919
$ws.Range("F7").Value = 1.870837887007969
$ws.Range("G7").Value = 46

$ws.Range("B8").Value = 0.1263763876690497
$ws.Range("C8").Value = 1.411360266901308
$ws.Range("D8").Value = 4.034073368543205
$ws.Range("E8").Value = 2.00850027845236
$ws.Range("F8").Value = 2.027171153542608
$ws.Range("G8").Value = 45

$ws.Range("B9").Value = 0.4927478854388472
$ws.Range("C9").Value = 1.400784718790886
$ws.Range("D9").Value = 4.048194332731589
$ws.Range("E9").Value = 2.012012508095213
$ws.Range("F9").Value = 1.973294551155114
$ws.Range("G9").Value = 44

# --- Row 10: update values; F10 did not previously exist, so it is newly populated ---
$ws.Range("B10").Value = 0.2097999035901463
$ws.Range("C10").Value = 1.453401524484553
$ws.Range("D10").Value = 3.844298923870199
$ws.Range("E10").Value = 1.960688380102815
$ws.Range("F10").Value = 1.97250243373573
$ws.Range("G10").Value = 43

# --- New row 11 for the Q9 forecast-error quantile ---
$ws.Range("A11").Value = "Q9"
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B11").Value = 0.6126507335481286
$ws.Range("C11").Value = 1.379826743076134
$ws.Range("D11").Value = 3.871268613803817
$ws.Range("E11").Value = 1.967553967189672
$ws.Range("F11").Value = 1.892404349385141
$ws.Range("G11").Value = 42
